$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new student record on row 22 ---
# Values are written in the same order the source application inserted them
# so that the shared-string table is built in the same sequence as the
# target workbook (name, class, fname, surname, lname, village, aadhar,
# email, gender, dob, photo path, gr_no).
$ws.Range("B22").Value = "Bhavani"
$ws.Range("C22").Value = "bcom"
$ws.Range("D22").Value = "Bhavani "
$ws.Range("F22").Value = "Gohil"
$ws.Range("E22").Value = "J"
$ws.Range("G22").Value = "chokdi"
$ws.Range("H22").Value = 8347078272
$ws.Range("I22").Value = "XXXXXXXX3584"
$ws.Range("J22").Value = "bhavanising008@gmail.com"

# Turn the e-mail address into a live hyperlink (this also introduces the
# "Hyperlink" cell style / underlined theme-color font used by the target
# workbook).
$ws.Hyperlinks.Add($ws.Range("J22"), "mailto:bhavanising008@gmail.com")

$ws.Range("K22").Value = "Male"

# DOB stored as a real date (serial 39646 == 2008-07-17) formatted with the
# built-in short-date number format (numFmtId 14).
$ws.Range("L22").Value = 39646
$ws.Range("L22").NumberFormat = "mm-dd-yy"

$ws.Range("M22").Value = "photos/2021-013.jpg"
$ws.Range("A22").Value = "2021-013"

# --- Update the view so the new row is visible/selected ---
$ws.Range("A22").Select() | Out-Null
